# #5: property boat&car done
# Fix up the "汽車" (Car) sheet: row 1 was accidentally populated with a
# copy of row 2's data instead of real column headers. Give it proper
# headers (matching the other property sheets) and add the extra
# property_category / category / date / legislator_name / legislator_id /
# source_file / index columns that the other sheets already have, plus a
# new "capacity" column (engine displacement) in place of "area".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# ---- Header row (row 1) ----------------------------------------------
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "capacity"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "register_date"
$ws.Cells.Item(1,6).Value = "register_reason"
$ws.Cells.Item(1,7).Value = "acquire_value"
$ws.Cells.Item(1,8).Value = "property_category"
$ws.Cells.Item(1,9).Value = "category"
$ws.Cells.Item(1,10).Value = "date"
$ws.Cells.Item(1,11).Value = "legislator_name"
$ws.Cells.Item(1,12).Value = "legislator_id"
$ws.Cells.Item(1,13).Value = "source_file"
$ws.Cells.Item(1,14).Value = "index"

# ---- Data row (row 2) new columns -------------------------------------
$ws.Cells.Item(2,8).Value = "land"
$ws.Cells.Item(2,9).Value = "normal"
$ws.Cells.Item(2,10).Value = "2013-12-30"
$ws.Cells.Item(2,11).Value = "管碧玲"
$ws.Cells.Item(2,12).Value = 1374
$ws.Cells.Item(2,13).Value = "tmpb8981"
$ws.Cells.Item(2,14).Value = 65

# ---- Formatting: reuse the existing header / data cell styles ---------
# Header style lives on B1:G1 (style index 1 - bold, bordered, centered)
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

# Data-row style lives on B2:G2 (style index 2)
$ws.Range("B2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

$excel.CutCopyMode = $false
